# "Generate Report for Handoff"
# The 1bbacf06-3ffc-4215-b44d-39c6a6eca9d1 file's row is removed from every
# sheet (it is now considered handed off / no longer handed-back), and the
# still-present 075d7a73... row's status + timestamps move from
# "Handed back: in sync with en-US" to "Ready for handoff".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"
$ov.Range("D2").Value = "2016-38-19 10:38:12"

# Drop the 1bbacf06... row entirely.
$ov.Rows.Item(3).Delete()

# Row deletion only touches sheetData; hyperlinks collection needs to be
# rebuilt by hand (this host deletes/re-adds hyperlinks sheet-wide).
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7f983055ada23edc3bbdc6dcba8544a48464a646/e2e/075d7a73-6826-41de-8bfe-e6cf8e17f9ec.md", "", "", "075d7a73-6826-41de-8bfe-e6cf8e17f9ec.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("E2").Value = "2016-03-19 10:38:09"

$zh.Rows.Item(3).Delete()

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7f983055ada23edc3bbdc6dcba8544a48464a646/e2e/075d7a73-6826-41de-8bfe-e6cf8e17f9ec.md", "", "", "075d7a73-6826-41de-8bfe-e6cf8e17f9ec.md")
$zh.Hyperlinks.Add($zh.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/7f983055ada23edc3bbdc6dcba8544a48464a646/e2e/075d7a73-6826-41de-8bfe-e6cf8e17f9ec.md", "", "", ".md")
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0da20fdcb2095ad1d9efce3b662a47508bef0b76/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/075d7a73-6826-41de-8bfe-e6cf8e17f9ec.94eadd793ec9569dbea66e55c42f4ff2303055cc.zh-cn.xlf", "", "", "075d7a73-6826-41de-8bfe-e6cf8e17f9ec.94eadd793ec9569dbea66e55c42f4ff2303055cc.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/f9d48f766536a39bd95aa2caf23c0530feaf8640/e2e/075d7a73-6826-41de-8bfe-e6cf8e17f9ec.md", "", "", "075d7a73-6826-41de-8bfe-e6cf8e17f9ec.md")
$zh.Hyperlinks.Add($zh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/580a580d4fa62df7e11ef69dcbb94f2df72a474a/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/075d7a73-6826-41de-8bfe-e6cf8e17f9ec.94eadd793ec9569dbea66e55c42f4ff2303055cc.zh-cn.xlf", "", "", "075d7a73-6826-41de-8bfe-e6cf8e17f9ec.94eadd793ec9569dbea66e55c42f4ff2303055cc.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = "Ready for handoff"
$de.Range("E2").Value = "2016-03-19 10:38:12"

$de.Rows.Item(3).Delete()

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7f983055ada23edc3bbdc6dcba8544a48464a646/e2e/075d7a73-6826-41de-8bfe-e6cf8e17f9ec.md", "", "", "075d7a73-6826-41de-8bfe-e6cf8e17f9ec.md")
$de.Hyperlinks.Add($de.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/7f983055ada23edc3bbdc6dcba8544a48464a646/e2e/075d7a73-6826-41de-8bfe-e6cf8e17f9ec.md", "", "", ".md")
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b037c97e912ee6d80f7a96fda6f39ce3f8dc28d4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/075d7a73-6826-41de-8bfe-e6cf8e17f9ec.94eadd793ec9569dbea66e55c42f4ff2303055cc.de-de.xlf", "", "", "075d7a73-6826-41de-8bfe-e6cf8e17f9ec.94eadd793ec9569dbea66e55c42f4ff2303055cc.de-de.xlf")
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/3715e864f0ffc4408a68ce04248c2c8760d6980f/e2e/075d7a73-6826-41de-8bfe-e6cf8e17f9ec.md", "", "", "075d7a73-6826-41de-8bfe-e6cf8e17f9ec.md")
$de.Hyperlinks.Add($de.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b9899ce00336e00ec3646ca43ba553cd838f7f7e/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/075d7a73-6826-41de-8bfe-e6cf8e17f9ec.94eadd793ec9569dbea66e55c42f4ff2303055cc.de-de.xlf", "", "", "075d7a73-6826-41de-8bfe-e6cf8e17f9ec.94eadd793ec9569dbea66e55c42f4ff2303055cc.de-de.xlf")

Write-Output "Done"
